# Update party name headers (row 1, columns B-P) to include full descriptive names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "IP/SJFL - Independence Party (Sjálfstæðisflokkur, IP/SJFL)"
$ws.Range("C1").Value = "PA - People's Alliance (Althydubandalag, PA)"
$ws.Range("D1").Value = "PP - Progressive Party (Framsóknarflokkur, PP)"
$ws.Range("E1").Value = "SDP - Social Democrats (Althqduflokkur, SDP)"
$ws.Range("F1").Value = "WA - Women's Alliance (Samtok um kvennalista, WA)"
$ws.Range("G1").Value = "PM - People's Movement (Thjodvaki - hreyfing f6lksins, PM)"
$ws.Range("H1").Value = "LG - Left-Greens (Vinstri græn, LG)"
$ws.Range("I1").Value = "LP - Liberal Party (Frjálslyndi flokkurinn, LP)"
$ws.Range("J1").Value = "SDA - Social Democratic Alliance (Samfylking, SDA)"
$ws.Range("K1").Value = "CM - Civic Movement (Borgarahreyfingin, CM)"
$ws.Range("L1").Value = "BF - Bright Future (Bjartar framtíðar, BF)"
$ws.Range("M1").Value = "P - Pirates (Pirata, P)"
$ws.Range("N1").Value = "VID - Reform (Viðreisn, VIÐ)"
$ws.Range("O1").Value = "M - The Central Party (Miðflokkurinn, M)"
$ws.Range("P1").Value = "RGPP - Right-Green People's Party (Hægri grænna, flokks fólksins, RGPP)"

# Clean up floating-point noise in row 10 (last data row) caused by re-computed
# minimal seat counts; these should be clean integers.
$ws.Range("B10").Value = 7
$ws.Range("D10").Value = 4
$ws.Range("H10").Value = 5
$ws.Range("M10").Value = 3
$ws.Range("O10").Value = 3
